# Bump the "Förändrad" (Changed) date in column C by one day
# for every data row (rows 2 through 451), changing serial date
# 45180 (2023-09-11) to 45181 (2023-09-12).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 451
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)  # Column C
    if ($cell.Value2() -eq 45180) {
        $cell.Value = 45181
    }
}
